$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 862.7917
$ws.Cells.Item(33, 9).Value = 885.3889
$ws.Cells.Item(33, 10).Value = 795
$ws.Cells.Item(33, 11).Value = 885.3889
$ws.Cells.Item(33, 12).Value = 795
$ws.Cells.Item(33, 13).Value = -656.3889
$ws.Cells.Item(33, 14).Value = -1253
$ws.Cells.Item(69, 8).Value = 476666.66
$ws.Cells.Item(69, 9).Value = 15000
$ws.Cells.Item(69, 10).Value = 569000
$ws.Cells.Item(69, 11).Value = 45000
$ws.Cells.Item(69, 12).Value = 1707000
$ws.Cells.Item(69, 13).Value = -44126
$ws.Cells.Item(69, 14).Value = -1708748
$ws.Cells.Item(72, 8).Value = 476666.66
$ws.Cells.Item(72, 9).Value = 15000
$ws.Cells.Item(72, 10).Value = 569000
$ws.Cells.Item(72, 11).Value = 135000
$ws.Cells.Item(72, 12).Value = 5121000
$ws.Cells.Item(72, 13).Value = -130632
$ws.Cells.Item(72, 14).Value = -5129736
$ws.Cells.Item(100, 8).Value = 80710.875
$ws.Cells.Item(100, 9).Value = 105000.836
$ws.Cells.Item(100, 10).Value = 7841
$ws.Cells.Item(100, 11).Value = 105000.836
$ws.Cells.Item(100, 12).Value = 7841
$ws.Cells.Item(100, 13).Value = -104459.836
$ws.Cells.Item(100, 14).Value = -8923
$ws.Cells.Item(121, 8).Value = 1845.6666
$ws.Cells.Item(121, 10).Value = 1845.6666
$ws.Cells.Item(121, 12).Value = 5536.9998
$ws.Cells.Item(121, 14).Value = -9030.9998
$ws.Cells.Item(137, 8).Value = 13595.214
$ws.Cells.Item(137, 9).Value = 39802.5
$ws.Cells.Item(137, 10).Value = 3112.3
$ws.Cells.Item(137, 11).Value = 119407.5
$ws.Cells.Item(137, 12).Value = 9336.900000000001
$ws.Cells.Item(137, 13).Value = -116857.5
$ws.Cells.Item(137, 14).Value = -14436.9
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2888.3704
$ws.Cells.Item(2, 9).Value = 2988.389
$ws.Cells.Item(2, 10).Value = 2688.3333
$ws.Cells.Item(2, 11).Value = 2988.389
$ws.Cells.Item(2, 12).Value = 2688.3333
$ws.Cells.Item(2, 13).Value = -2875.389
$ws.Cells.Item(2, 14).Value = -2914.3333
$ws.Cells.Item(32, 8).Value = 20230.625
$ws.Cells.Item(32, 9).Value = 21171.264
$ws.Cells.Item(32, 11).Value = 21171.264
$ws.Cells.Item(32, 13).Value = -20884.264
$ws.Cells.Item(43, 8).Value = 69914.664
$ws.Cells.Item(43, 10).Value = 92372
$ws.Cells.Item(43, 12).Value = 92372
$ws.Cells.Item(43, 14).Value = -92998
$ws.Cells.Item(45, 8).Value = 2886.92
$ws.Cells.Item(45, 9).Value = 2083
$ws.Cells.Item(45, 10).Value = 4092.8
$ws.Cells.Item(45, 11).Value = 2083
$ws.Cells.Item(45, 12).Value = 4092.8
$ws.Cells.Item(45, 13).Value = -1706
$ws.Cells.Item(45, 14).Value = -4846.8
$ws.Cells.Item(61, 8).Value = 2466.6667
$ws.Cells.Item(61, 9).Value = 1262.0714
$ws.Cells.Item(61, 11).Value = 1262.0714
$ws.Cells.Item(61, 13).Value = -1050.0714
$ws.Cells.Item(97, 8).Value = 1834.2222
$ws.Cells.Item(97, 9).Value = 1111.2609
$ws.Cells.Item(97, 11).Value = 1111.2609
$ws.Cells.Item(97, 13).Value = -615.2609
$ws.Cells.Item(116, 8).Value = 2888.3704
$ws.Cells.Item(116, 9).Value = 2988.389
$ws.Cells.Item(116, 10).Value = 2688.3333
$ws.Cells.Item(116, 11).Value = 2988.389
$ws.Cells.Item(116, 12).Value = 2688.3333
$ws.Cells.Item(116, 13).Value = -694.3890000000001
$ws.Cells.Item(116, 14).Value = -7276.3333
$ws.Cells.Item(132, 8).Value = 1838.8334
$ws.Cells.Item(132, 9).Value = 1533.7858
$ws.Cells.Item(132, 10).Value = 2906.5
$ws.Cells.Item(132, 11).Value = 4601.357400000001
$ws.Cells.Item(132, 12).Value = 8719.5
$ws.Cells.Item(132, 13).Value = -2071.357400000001
$ws.Cells.Item(132, 14).Value = -13779.5
$ws.Cells.Item(136, 8).Value = 2466.6667
$ws.Cells.Item(136, 9).Value = 1262.0714
$ws.Cells.Item(136, 11).Value = 3786.2142
$ws.Cells.Item(136, 13).Value = -1236.2142
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2888.3704
$ws.Cells.Item(3, 9).Value = 2988.389
$ws.Cells.Item(3, 10).Value = 2688.3333
$ws.Cells.Item(3, 11).Value = 2988.389
$ws.Cells.Item(3, 12).Value = 2688.3333
$ws.Cells.Item(3, 13).Value = -2874.389
$ws.Cells.Item(3, 14).Value = -2916.3333
$ws.Cells.Item(75, 8).Value = 30713.857
$ws.Cells.Item(75, 9).Value = 19998.5
$ws.Cells.Item(75, 11).Value = 19998.5
$ws.Cells.Item(75, 13).Value = -19062.5
$ws.Cells.Item(78, 8).Value = 30713.857
$ws.Cells.Item(78, 9).Value = 19998.5
$ws.Cells.Item(78, 11).Value = 59995.5
$ws.Cells.Item(78, 13).Value = -55315.5
$ws.Cells.Item(94, 8).Value = 649.125
$ws.Cells.Item(94, 9).Value = 551.3182
$ws.Cells.Item(94, 10).Value = 1725
$ws.Cells.Item(94, 11).Value = 551.3182
$ws.Cells.Item(94, 12).Value = 1725
$ws.Cells.Item(94, 13).Value = -100.3182
$ws.Cells.Item(94, 14).Value = -2627
$ws.Cells.Item(99, 8).Value = 1914.1818
$ws.Cells.Item(99, 9).Value = 1486.7858
$ws.Cells.Item(99, 10).Value = 2662.125
$ws.Cells.Item(99, 11).Value = 1486.7858
$ws.Cells.Item(99, 12).Value = 2662.125
$ws.Cells.Item(99, 13).Value = 11.21419999999989
$ws.Cells.Item(99, 14).Value = -5658.125
$ws.Cells.Item(140, 8).Value = 99931.53
$ws.Cells.Item(140, 10).Value = 99931.53
$ws.Cells.Item(140, 12).Value = 99931.53
$ws.Cells.Item(140, 14).Value = -110291.53
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1368
$ws.Cells.Item(16, 9).Value = 1332
$ws.Cells.Item(16, 10).Value = 1476
$ws.Cells.Item(16, 11).Value = 1332
$ws.Cells.Item(16, 12).Value = 1476
$ws.Cells.Item(16, 13).Value = -1045
$ws.Cells.Item(16, 14).Value = -2050
$ws.Cells.Item(94, 8).Value = 1235
$ws.Cells.Item(94, 9).Value = 1193.75
$ws.Cells.Item(94, 10).Value = 1400
$ws.Cells.Item(94, 11).Value = 1193.75
$ws.Cells.Item(94, 12).Value = 1400
$ws.Cells.Item(94, 13).Value = -742.75
$ws.Cells.Item(94, 14).Value = -2302
$ws.Cells.Item(113, 8).Value = 1368
$ws.Cells.Item(113, 9).Value = 1332
$ws.Cells.Item(113, 10).Value = 1476
$ws.Cells.Item(113, 11).Value = 1332
$ws.Cells.Item(113, 12).Value = 1476
$ws.Cells.Item(113, 13).Value = 838
$ws.Cells.Item(113, 14).Value = -5816
$ws.Cells.Item(122, 8).Value = 2000
$ws.Cells.Item(122, 9).Value = 2000
$ws.Cells.Item(122, 11).Value = 6000
$ws.Cells.Item(122, 13).Value = -3550
$ws.Cells.Item(134, 8).Value = 2941
$ws.Cells.Item(134, 9).Value = 2470.75
$ws.Cells.Item(134, 10).Value = 3646.375
$ws.Cells.Item(134, 11).Value = 7412.25
$ws.Cells.Item(134, 12).Value = 10939.125
$ws.Cells.Item(134, 13).Value = -4877.25
$ws.Cells.Item(134, 14).Value = -16009.125
$ws.Cells.Item(141, 8).Value = 175372.61
$ws.Cells.Item(141, 10).Value = 175372.61
$ws.Cells.Item(141, 12).Value = 175372.61
$ws.Cells.Item(141, 14).Value = -185732.61
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 221.72728
$ws.Cells.Item(2, 9).Value = 106.4
$ws.Cells.Item(2, 11).Value = 638.4000000000001
$ws.Cells.Item(2, 13).Value = -525.4000000000001
$ws.Cells.Item(81, 8).Value = 1816.3334
$ws.Cells.Item(81, 9).Value = 1874.75
$ws.Cells.Item(81, 10).Value = 1699.5
$ws.Cells.Item(81, 11).Value = 5624.25
$ws.Cells.Item(81, 12).Value = 5098.5
$ws.Cells.Item(81, 13).Value = -4501.25
$ws.Cells.Item(81, 14).Value = -7344.5
$ws.Cells.Item(84, 8).Value = 1816.3334
$ws.Cells.Item(84, 9).Value = 1874.75
$ws.Cells.Item(84, 10).Value = 1699.5
$ws.Cells.Item(84, 11).Value = 16872.75
$ws.Cells.Item(84, 12).Value = 15295.5
$ws.Cells.Item(84, 13).Value = -11256.75
$ws.Cells.Item(84, 14).Value = -26527.5
$ws.Cells.Item(114, 8).Value = 666.6923
$ws.Cells.Item(114, 9).Value = 644.6
$ws.Cells.Item(114, 11).Value = 1933.8
$ws.Cells.Item(114, 13).Value = 1320.2
$ws.Cells.Item(131, 8).Value = 1892.963
$ws.Cells.Item(131, 9).Value = 1745.4286
$ws.Cells.Item(131, 11).Value = 5236.2858
$ws.Cells.Item(131, 13).Value = -196.2857999999997
$ws.Cells.Item(132, 8).Value = 2047.4615
$ws.Cells.Item(132, 9).Value = 5054
$ws.Cells.Item(132, 10).Value = 1500.8182
$ws.Cells.Item(132, 11).Value = 45486
$ws.Cells.Item(132, 12).Value = 13507.3638
$ws.Cells.Item(132, 13).Value = -42956
$ws.Cells.Item(132, 14).Value = -18567.3638
$ws.Cells.Item(133, 8).Value = 5676.3335
$ws.Cells.Item(133, 9).Value = 5676.3335
$ws.Cells.Item(133, 11).Value = 17029.0005
$ws.Cells.Item(133, 13).Value = -11969.0005
$ws.Cells.Item(139, 8).Value = 6273.5386
$ws.Cells.Item(139, 10).Value = 2500
$ws.Cells.Item(139, 12).Value = 7500
$ws.Cells.Item(139, 14).Value = -17780
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 2305.3333
$ws.Cells.Item(31, 9).Value = 1374.75
$ws.Cells.Item(31, 11).Value = 1374.75
$ws.Cells.Item(31, 13).Value = -1082.75
$ws.Cells.Item(37, 8).Value = 2305.3333
$ws.Cells.Item(37, 9).Value = 1374.75
$ws.Cells.Item(37, 11).Value = 1374.75
$ws.Cells.Item(37, 13).Value = -1097.75
$ws.Cells.Item(97, 8).Value = 1332.9259
$ws.Cells.Item(97, 10).Value = 1822.625
$ws.Cells.Item(97, 12).Value = 1822.625
$ws.Cells.Item(97, 14).Value = -2814.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 5939.4287
$ws.Cells.Item(35, 9).Value = 5426.4165
$ws.Cells.Item(35, 10).Value = 9017.5
$ws.Cells.Item(35, 11).Value = 5426.4165
$ws.Cells.Item(35, 12).Value = 9017.5
$ws.Cells.Item(35, 13).Value = -5090.4165
$ws.Cells.Item(35, 14).Value = -9689.5
$ws.Cells.Item(46, 8).Value = 6334
$ws.Cells.Item(46, 9).Value = 2956.8572
$ws.Cells.Item(46, 10).Value = 8698
$ws.Cells.Item(46, 11).Value = 2956.8572
$ws.Cells.Item(46, 12).Value = 8698
$ws.Cells.Item(46, 13).Value = -2768.8572
$ws.Cells.Item(46, 14).Value = -9074
$ws.Cells.Item(108, 8).Value = 72949.5
$ws.Cells.Item(108, 10).Value = 72949.5
$ws.Cells.Item(108, 12).Value = 72949.5
$ws.Cells.Item(108, 14).Value = -80629.5
$ws.Cells.Item(122, 8).Value = 5297.222
$ws.Cells.Item(122, 9).Value = 4556.533
$ws.Cells.Item(122, 11).Value = 13669.599
$ws.Cells.Item(122, 13).Value = -11219.599
$ws.Cells.Item(132, 8).Value = 4974
$ws.Cells.Item(132, 9).Value = 1999.5
$ws.Cells.Item(132, 11).Value = 5998.5
$ws.Cells.Item(132, 13).Value = -3468.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1976.5
$ws.Cells.Item(96, 10).Value = 2610.75
$ws.Cells.Item(96, 12).Value = 2610.75
$ws.Cells.Item(96, 14).Value = -5356.75
$ws.Cells.Item(119, 8).Value = 34600
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()
$ws.Cells.Item(124, 8).Value = 50429
$ws.Cells.Item(124, 10).Value = 50429
$ws.Cells.Item(124, 12).Value = 50429
$ws.Cells.Item(124, 14).Value = -60249
$ws.Cells.Item(132, 8).Value = 18745.309
$ws.Cells.Item(132, 9).Value = 21716
$ws.Cells.Item(132, 10).Value = 921.1667
$ws.Cells.Item(132, 11).Value = 65148
$ws.Cells.Item(132, 12).Value = 2763.5001
$ws.Cells.Item(132, 13).Value = -62618
$ws.Cells.Item(132, 14).Value = -7823.5001
$ws.Cells.Item(136, 8).Value = 20728.162
$ws.Cells.Item(136, 9).Value = 24507.064
$ws.Cells.Item(136, 10).Value = 1203.8334
$ws.Cells.Item(136, 11).Value = 73521.192
$ws.Cells.Item(136, 12).Value = 3611.5002
$ws.Cells.Item(136, 13).Value = -70971.192
